$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Sheet1" to "Bus Load"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Bus Load"

# Add a new worksheet named "TestSheet" after the first sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TestSheet"

# Put content in the new sheet
$ws2.Range("A1").Value = "Tstt"

# Make the new sheet the active one
$ws2.Activate()
